$wb = $excel.ActiveWorkbook

# --- Base sheet we are branching the two new modules from ---
$base = $wb.Worksheets.Item("GroupLoanApplication")

# --- Create "GroupLoanAgreement" as a copy placed right after GroupLoanApplication ---
$base.Copy([System.Reflection.Missing]::Value, $base)
$agreement = $wb.Worksheets.Item($base.Index + 1)
$agreement.Name = "GroupLoanAgreement"
# Drop the 5th (Amount) column that GroupLoanApplication had - the new module only needs A:D
$agreement.Columns.Item(5).Delete()
$agreement.Range("A2").Value = "Group Loan Agreement"

# --- Create "GroupLoanSanction" as a copy placed right after GroupLoanAgreement ---
$agreement.Copy([System.Reflection.Missing]::Value, $agreement)
$sanction = $wb.Worksheets.Item($agreement.Index + 1)
$sanction.Name = "GroupLoanSanction"
$sanction.Range("A2").Value = "Group Loan Sanction"
$sanction.Rows.Item(2).RowHeight = 45

# --- Restore / adjust the selection shown on each of the three module sheets ---
$base.Activate()
$base.Range("A1:D2").Select()

$agreement.Activate()
$agreement.Range("A2").Select()

$sanction.Activate()
$sanction.Range("L9").Select()
